# Updated cryptos list on Tue Jul 11 19:25:27 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns scraped from coinranking.com,
# and fixes the Litecoin / WrappedEther row ordering (rows 13-14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.531.43'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.873.01'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'247.54"
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'0.4736"
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('D8').Value = "'0.2894"
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = "'0.06467"
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').Value = "'22.00"
$ws.Range('E10').Value = '  +2.74%  '
$ws.Range('D11').Value = "'0.07718"
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = "'0.7389"
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = "'96.05"
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.871.53'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = "'5.166"
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').Value = "'274.65"
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '30.609.25'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = "'13.24"
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = "'0.000007472"
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '2.112.00'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').Value = "'0.9995"
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = "'5.214"
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = "'6.160"
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').Value = "'165.44"
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = "'9.182"
$ws.Range('E26').Value = '  -1.52%  '
$ws.Range('D27').Value = "'18.66"
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').Value = "'1.900"
$ws.Range('E28').Value = '  -4.92%  '
$ws.Range('D29').Value = "'0.09954"
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('D30').Value = "'1.344"
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').Value = "'1.506"
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').Value = "'4.233"
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('D33').Value = "'4.076"
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').Value = "'0.04761"
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').Value = "'1.118"
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('D36').Value = "'0.6913"
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').Value = "'2.715"
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('D39').Value = "'2.753"
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').Value = "'6.244"
$ws.Range('D41').Value = "'73.27"
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('D42').Value = "'1.966"
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = "'1.000"
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = "'0.4152"
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').Value = "'0.8336"
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('D46').Value = "'101.28"
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').Value = "'9.340"
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('D48').Value = "'35.27"
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = "'6.962"
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('D50').Value = "'915.36"
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('E51').Value = '  +0.99%  '

# The leading apostrophe above also stamps a quote-prefix cell style; reset
# each of those cells back to the (unstyled) Normal style to match the source.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
